# Automated scheduled-task update for datos_dropcontrol/2025-07-28.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3, column A: corrected timestamp (sub-millisecond precision fix)
$ws.Range("A3").Value = 45866.08356709491

# New row 4: next sensor reading appended by the scheduled task
$ws.Range("A4").NumberFormat = $ws.Range("A3").NumberFormat()
$ws.Range("A4").Value = 45866.16688976369
$ws.Range("B4").Value = 2025
$ws.Range("C4").Value = 31
$ws.Range("D4").Value = 13.61
$ws.Range("E4").Value = 91.53
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 2.89
$ws.Range("H4").Value = "ESE"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "04:00:19"
